$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (test case template columns) ---
$ws.Range("A1").Value = "Step"
$ws.Range("B1").Value = "Action to do"
$ws.Range("C1").Value = "Expected result"
$ws.Range("D1").Value = "Actual result"
$ws.Range("E1").Value = "Test result"

# --- First test-case row ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "SELECT flightNumber FROM flight WHERE departureAirport='London' AND arrivalAirport='Munich' AND averageTicketPrice<100 AND availableSeats>4"
$ws.Range("C2").Value = "1001;2001;3001"

# Wrap the long SQL text and grow the row to show it
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 90

# --- Column widths to fit the new content ---
$ws.Columns.Item(1).ColumnWidth = 4.17
$ws.Columns.Item(2).ColumnWidth = 29.33
$ws.Columns.Item(3).ColumnWidth = 18.66
$ws.Columns.Item(4).ColumnWidth = 17.17
$ws.Columns.Item(5).ColumnWidth = 21.33

# --- Page setup for printing ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave selection on D2, like after tabbing past the filled cells ---
$null = $ws.Range("D2").Select()
